# Generate Report for Handback
#
# This script reflects a localization "handback" run: files that were
# previously only "handed off" for translation now have a completed
# handback recorded (Latest Target File / Latest Handback File /
# Latest Handback DateTime columns populated), and the Overview status
# message is updated to reflect the new state.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: refresh the status text for both language columns ---
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------
$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

$zhHandoffXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhTargetUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/032df44b9e6978247781691418108b8f030417bc/e2e/a.md"
$zhHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9da02e668c5a9683e25e90472a4307ad1e1b1904/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/hb/$zhHandoffXlf"

# Row 2 (a.md): record the resulting target file + handback file + datetime
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), $zhTargetUrl, $null, $null, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), $zhHandbackUrl, $null, $null, $zhHandoffXlf)
$zhcn.Range("G2").Value = "2016-02-18 03:34:40"

# Row 3 (b.md): same handback batch
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), $zhTargetUrl, $null, $null, "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), $zhHandbackUrl, $null, $null, $zhHandoffXlf)
$zhcn.Range("G3").Value = "2016-02-18 03:34:40"

# --- de-de sheet -------------------------------------------------------
$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

$deHandoffXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deTargetUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/032df44b9e6978247781691418108b8f030417bc/e2e/a.md"
$deHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/41281aea5de1421f3c8cff85e3b96bb26a49a93c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/hb/$deHandoffXlf"

# Row 2 (a.md)
$dede.Hyperlinks.Add($dede.Range("E2"), $deTargetUrl, $null, $null, "a.md")
$dede.Hyperlinks.Add($dede.Range("F2"), $deHandbackUrl, $null, $null, $deHandoffXlf)
$dede.Range("G2").Value = "2016-02-18 03:35:00"

# Row 3 (b.md)
$dede.Hyperlinks.Add($dede.Range("E3"), $deTargetUrl, $null, $null, "a.md")
$dede.Hyperlinks.Add($dede.Range("F3"), $deHandbackUrl, $null, $null, $deHandoffXlf)
$dede.Range("G3").Value = "2016-02-18 03:35:00"

Write-Output "Handback report generated."
